$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (serial 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = 45436

# Update price values in column D (rows 33-37)
$ws.Range("D33").Value = 767.647
$ws.Range("D34").Value = 1139.001
$ws.Range("D35").Value = 1427.198
$ws.Range("D36").Value = 1718.204
$ws.Range("D37").Value = 1878.96
